$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing "Late"/"Waived"/"Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()

# The newly inserted column keeps the width Excel computed for it (11
# characters), but unlike its bestFit neighbours it is an explicit,
# non-autofit width.
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Update the selection to match where the user left the cursor after the
# edit, then make this sheet the active / selected tab.
$ws.Range("Q11").Select()
$ws.Activate()
